# Update the icon artwork on slide 3:
#  - shift/round the rounded-square backing shape
#  - reposition + resize the icon group to match the new artwork

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Shape 1: "圆角矩形 12" (rounded-square logo background)
$backing = $s.Shapes.Item(1)
$backing.Left = 219.40001
$backing.Adjustments.Item(1) = 0.24735

# Shape 2: "组合 11" (grouped icon artwork)
$icon = $s.Shapes.Item(2)
$icon.Left = 325.85
$icon.Top = 112.15
$icon.Width = 288.80001
$icon.Height = 315.7
